$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 526, shifting existing rows 526:551 down to 527:552
$ws.Rows("526:526").Insert()

# Populate the newly inserted row 526 with the new record
$ws.Cells.Item(526, 1).Value = 4
$ws.Cells.Item(526, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(526, 3).Value = "Los Lagos"
$ws.Cells.Item(526, 4).Value = 45041
$ws.Cells.Item(526, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(526, 5).Value = 10
$ws.Cells.Item(526, 6).Value = 100114013
$ws.Cells.Item(526, 7).Value = "Zanahoria"
$ws.Cells.Item(526, 8).Value = "Sin especificar"
$ws.Cells.Item(526, 9).Value = "Primera"
$ws.Cells.Item(526, 10).Value = 900
$ws.Cells.Item(526, 11).Value = 8000
$ws.Cells.Item(526, 12).Value = 8500
$ws.Cells.Item(526, 13).Value = 8250
$ws.Cells.Item(526, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(526, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(526, 16).Value = 412
$ws.Cells.Item(526, 17).Value = 20
$ws.Cells.Item(526, 18).Value = "Hortaliza"
